# Update "展览" (exhibition) and "全部类型" (all types) sheets with refreshed
# attendance/price numbers and a newly-announced event
# (合肥·第十四届次元之门动漫游戏博览会) that slots in on 2024-07-20.

$wb = $excel.ActiveWorkbook

# Helper: assign a text value to a cell while avoiding Excel's automatic
# number/date sniffing (e.g. "2024-07-20" being turned into a date serial),
# and without leaving a stray explicit NumberFormat-driven style behind.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(2, 6).Value = 68
$ws1.Cells.Item(5, 6).Value = 24
$ws1.Cells.Item(6, 6).Value = 18
$ws1.Cells.Item(8, 6).Value = 7937
$ws1.Cells.Item(9, 6).Value = 752
$ws1.Cells.Item(10, 6).Value = 225
$ws1.Cells.Item(12, 6).Value = 758
$ws1.Cells.Item(15, 6).Value = 199

# Insert the new row 16 (pushing 赛马娘/MAX特摄/环形宇宙 down by one),
# carrying the bold/centered/bordered style used by column A down too.
$ws1.Rows.Item(16).Insert()
$ws1.Cells.Item(15, 1).Copy()
$ws1.Cells.Item(16, 1).PasteSpecial(-4122)

$ws1.Cells.Item(16, 1).Value = 15
Set-TextValue $ws1.Cells.Item(16, 2) "2024-07-20"
$ws1.Cells.Item(16, 3).Value = "合肥·第十四届次元之门动漫游戏博览会"
$ws1.Cells.Item(16, 4).Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Cells.Item(16, 5).Value = "2024.07.20 10:00-07.21 17:00"
$ws1.Cells.Item(16, 6).Value = 14
$ws1.Cells.Item(16, 7).Value = 68
$ws1.Cells.Item(16, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85336"
$ws1.Cells.Item(16, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/Bu6iQPJ01715161445356.jpeg"

# Rows below the insertion point shifted down by one; renumber the
# leading index column and update their refreshed "想去人数" counts at
# their new row positions.
$ws1.Cells.Item(17, 1).Value = 16
$ws1.Cells.Item(17, 6).Value = 41
$ws1.Cells.Item(18, 1).Value = 17
$ws1.Cells.Item(19, 1).Value = 18
$ws1.Cells.Item(19, 6).Value = 822

# ---------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(2, 6).Value = 68
$ws4.Cells.Item(5, 6).Value = 24
$ws4.Cells.Item(6, 6).Value = 18
$ws4.Cells.Item(9, 6).Value = 7937
$ws4.Cells.Item(10, 6).Value = 752
$ws4.Cells.Item(11, 6).Value = 225
$ws4.Cells.Item(13, 6).Value = 758
$ws4.Cells.Item(16, 6).Value = 199

# Insert the new row 17 (pushing 赛马娘/MAX特摄/环形宇宙/菊次郎的夏天 down
# by one), carrying the column-A style down too.
$ws4.Rows.Item(17).Insert()
$ws4.Cells.Item(16, 1).Copy()
$ws4.Cells.Item(17, 1).PasteSpecial(-4122)

$ws4.Cells.Item(17, 1).Value = 16
Set-TextValue $ws4.Cells.Item(17, 2) "2024-07-20"
$ws4.Cells.Item(17, 3).Value = "合肥·第十四届次元之门动漫游戏博览会"
$ws4.Cells.Item(17, 4).Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Cells.Item(17, 5).Value = "2024.07.20 10:00-07.21 17:00"
$ws4.Cells.Item(17, 6).Value = 14
$ws4.Cells.Item(17, 7).Value = 68
$ws4.Cells.Item(17, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85336"
$ws4.Cells.Item(17, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/Bu6iQPJ01715161445356.jpeg"

# Rows below the insertion point shifted down by one; renumber the
# leading index column and update their refreshed "想去人数" counts at
# their new row positions.
$ws4.Cells.Item(18, 1).Value = 17
$ws4.Cells.Item(18, 6).Value = 41
$ws4.Cells.Item(19, 1).Value = 18
$ws4.Cells.Item(20, 1).Value = 19
$ws4.Cells.Item(20, 6).Value = 823
$ws4.Cells.Item(21, 1).Value = 20
